$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Rename tag variables to their "_NCI" suffixed versions per Benedikt Merz e-mail
$ws.Range("B20").Value = "sugars_NCI"
$ws.Range("B21").Value = "KMT_tag_NCI"
$ws.Range("B22").Value = "KMF_tag_NCI"

# Update the active selection to reflect where the edit was made
$ws.Activate()
$ws.Range("C24").Select()
